$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2
    $newVal = $oldVal - 20000
    $cell.Value = $newVal
}
